$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B4").Value = 120
$ws.Range("B5").Value = 0.125
$ws.Range("B24").Value = 240
$ws.Range("B26").Value = 3

# Add new rows 28-29
$ws.Range("A28").Value = "notchWidth"
$ws.Range("B28").Value = 10
$ws.Range("C28").Value = "mm"

$ws.Range("A29").Value = "spoolThickness"
$ws.Range("B29").Value = 0.22
$ws.Range("C29").Value = "in"

# Update selection to B6
$ws.Range("B6").Select()
